$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.913.85"
$ws.Range("E2").Value = "  -2.77%  "
$ws.Range("D3").Value = "1.858.75"
$ws.Range("E3").Value = "  -2.19%  "
$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'305.27"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").Value = "'0.9988"
$ws.Range("D7").Value = "'0.5042"
$ws.Range("E7").Value = "  -2.73%  "
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("D9").Value = "'0.07116"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").Value = "'0.8841"
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("D11").Value = "'20.50"
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07562"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.861.58"
$ws.Range("E13").Value = "  -2.15%  "
$ws.Range("D14").Value = "'5.279"
$ws.Range("E14").Value = "  -2.95%  "
$ws.Range("D15").Value = "'88.93"
$ws.Range("E15").Value = "  -3.36%  "
$ws.Range("D16").Value = "'0.9994"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "'0.000008344"
$ws.Range("E17").Value = "  -4.23%  "
$ws.Range("D18").Value = "'0.9985"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "'14.03"
$ws.Range("E19").Value = "  -2.94%  "
$ws.Range("D20").Value = "26.968.22"
$ws.Range("E20").Value = "  -2.70%  "
$ws.Range("D21").Value = "'5.017"
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("D22").Value = "2.110.08"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").Value = "'10.45"
$ws.Range("E23").Value = "  -3.42%  "
$ws.Range("D24").Value = "'6.450"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("D25").Value = "'1.846"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("D26").Value = "'146.63"
$ws.Range("E26").Value = "  -4.85%  "
$ws.Range("D27").Value = "'17.91"
$ws.Range("E27").Value = "  -2.14%  "
$ws.Range("E28").Value = "  -4.51%  "
$ws.Range("D29").Value = "'112.34"
$ws.Range("E29").Value = "  -2.09%  "
$ws.Range("D30").Value = "'4.642"
$ws.Range("E30").Value = "  -4.11%  "
$ws.Range("D31").Value = "'4.642"
$ws.Range("E31").Value = "  -3.15%  "
$ws.Range("D32").Value = "'0.09020"
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("D34").Value = "'3.039"
$ws.Range("E34").Value = "  -4.52%  "
$ws.Range("D35").Value = "'1.146"
$ws.Range("E35").Value = "  -7.06%  "
$ws.Range("D36").Value = "'0.7229"
$ws.Range("E36").Value = "  -6.97%  "
$ws.Range("D37").Value = "'0.02033"
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("D38").Value = "'3.034"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("D39").Value = "'2.449"
$ws.Range("E39").Value = "  -6.38%  "
$ws.Range("D40").Value = "'1.070"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("D41").Value = "'0.5252"
$ws.Range("E41").Value = "  -4.29%  "
$ws.Range("D42").Value = "'6.525"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("D43").Value = "'114.89"
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("D44").Value = "'8.211"
$ws.Range("E44").Value = "  -3.14%  "
$ws.Range("D45").Value = "'0.1459"
$ws.Range("E45").Value = "  -2.98%  "
$ws.Range("D46").Value = "'0.9983"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").Value = "'0.4575"
$ws.Range("E47").Value = "  -4.27%  "
$ws.Range("D48").Value = "'9.892"
$ws.Range("E48").Value = "  -5.09%  "
$ws.Range("D49").Value = "'1.550"
$ws.Range("E49").Value = "  -3.80%  "
$ws.Range("D50").Value = "'36.36"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "'63.81"
$ws.Range("E51").Value = "  -4.19%  "
